# Remove hard-coded column names.
# The "氏名" (Name) header used for the teacher-name column is replaced
# with "教員名" (Teacher Name) so the column label is no longer a
# hard-coded/ambiguous "name" string. This updates the shared-string
# table, the worksheet header cell, and (automatically, since the range
# is the header of the Excel Table) the table's column name as well.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "教員名"
